$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel, so they remain text like the rest
# of the Price column.
$textCells = @("D4", "D5", "D6", "D8", "D11", "D12", "D14", "D20", "D21", "D23", "D24", "D25", "D28", "D29", "D32", "D34", "D36", "D37", "D38", "D40", "D41", "D43", "D45", "D51")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '67.749.11'
$ws.Range("E2").Value = '  -1.09%  '
$ws.Range("D3").Value = '3.332.61'
$ws.Range("E3").Value = '  -1.05%  '
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '582.43'
$ws.Range("E5").Value = '  -1.88%  '
$ws.Range("D6").Value = '176.17'
$ws.Range("E6").Value = '  -5.35%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '0.589'
$ws.Range("E8").Value = '  -1.59%  '
$ws.Range("D9").Value = '3.328.09'
$ws.Range("E9").Value = '  -0.95%  '
$ws.Range("E10").Value = '  -3.03%  '
$ws.Range("D11").Value = '0.577'
$ws.Range("E11").Value = '  -1.60%  '
$ws.Range("D12").Value = '45.65'
$ws.Range("E12").Value = '  -3.50%  '
$ws.Range("E13").Value = '  -3.76%  '
$ws.Range("D14").Value = '663.38'
$ws.Range("E14").Value = '  +3.24%  '
$ws.Range("D15").Value = '3.871.97'
$ws.Range("E15").Value = '  -0.84%  '
$ws.Range("E16").Value = '  -1.67%  '
$ws.Range("D17").Value = '67.909.11'
$ws.Range("E17").Value = '  -0.94%  '
$ws.Range("E18").Value = '  -0.98%  '
$ws.Range("D19").Value = '3.333.59'
$ws.Range("E19").Value = '  -1.15%  '
$ws.Range("D20").Value = '17.47'
$ws.Range("E20").Value = '  -2.69%  '
$ws.Range("D21").Value = '10.96'
$ws.Range("E21").Value = '  -1.29%  '
$ws.Range("E22").Value = '  -2.28%  '
$ws.Range("D23").Value = '5.43'
$ws.Range("E23").Value = '  +6.92%  '
$ws.Range("D24").Value = '17.11'
$ws.Range("E24").Value = '  -4.68%  '
$ws.Range("D25").Value = '99.42'
$ws.Range("E25").Value = '  +0.00%  '
$ws.Range("E26").Value = '  -5.90%  '
$ws.Range("E27").Value = '  -5.97%  '
$ws.Range("D28").Value = '9.30'
$ws.Range("E28").Value = '  -4.83%  '
$ws.Range("D29").Value = '33.63'
$ws.Range("E29").Value = '  +2.21%  '
$ws.Range("E30").Value = '  +8.98%  '
$ws.Range("E31").Value = '  -2.77%  '
$ws.Range("D32").Value = '592.47'
$ws.Range("E32").Value = '  -3.26%  '
$ws.Range("E33").Value = '  -1.35%  '
$ws.Range("D34").Value = '0.105'
$ws.Range("E34").Value = '  -1.23%  '
$ws.Range("B35").Value = 'Maker'
$ws.Range("C35").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D35").Value = '3.726.81'
$ws.Range("E35").Value = '  -6.36%  '
$ws.Range("B36").Value = 'Dai'
$ws.Range("C36").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.15%  '
$ws.Range("D37").Value = '56.87'
$ws.Range("E37").Value = '  +1.09%  '
$ws.Range("D38").Value = '3.29'
$ws.Range("E38").Value = '  -11.28%  '
$ws.Range("E39").Value = '  +0.55%  '
$ws.Range("D40").Value = '33.59'
$ws.Range("E40").Value = '  -0.31%  '
$ws.Range("D41").Value = '2.64'
$ws.Range("E41").Value = '  -5.75%  '
$ws.Range("E42").Value = '  -6.23%  '
$ws.Range("D43").Value = '0.334'
$ws.Range("E43").Value = '  -2.71%  '
$ws.Range("D44").Value = '0.0₃0667'
$ws.Range("E44").Value = '  -5.47%  '
$ws.Range("D45").Value = '3.24'
$ws.Range("E45").Value = '  -4.96%  '
$ws.Range("E46").Value = '  -3.67%  '
$ws.Range("E47").Value = '  +0.14%  '
$ws.Range("E48").Value = '  -1.64%  '
$ws.Range("E49").Value = '  -0.19%  '
$ws.Range("E50").Value = '  -0.63%  '
$ws.Range("D51").Value = '127.14'
$ws.Range("E51").Value = '  -3.17%  '
